# This workbook contains a weekly price table for "Coliflor" at the
# "Vega Modelo de Temuco" market. Two new weekly observations were added
# to the data set. In the source data these new rows are inserted right
# after the existing row 625 (i.e. they become the new rows 626 and 627),
# pushing all the previously-existing rows (626-723) down by two rows
# (they become rows 628-725). The dimension grows from A1:R723 to
# A1:R725.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 626, shifting existing data (rows
# 626:723) down to rows 628:725.
$ws.Rows("626:627").Insert()

# --- New row 626 ---
$ws.Range("A626").Value = 10
$ws.Range("B626").Value = "Vega Modelo de Temuco"
$ws.Range("C626").Value = "La Araucanía"
$ws.Range("D626").Value = 45180
$ws.Range("E626").Value = 9
$ws.Range("F626").Value = 100112008
$ws.Range("G626").Value = "Coliflor"
$ws.Range("H626").Value = "Sin especificar"
$ws.Range("I626").Value = "Primera"
$ws.Range("J626").Value = 1200
$ws.Range("K626").Value = 1300
$ws.Range("L626").Value = 1300
$ws.Range("M626").Value = 1300
$ws.Range("N626").Value = "`$/unidad"
$ws.Range("O626").Value = "Provincia del Elquí"
$ws.Range("P626").Value = 1300
$ws.Range("Q626").Value = 1
$ws.Range("R626").Value = "Hortaliza"

# --- New row 627 ---
$ws.Range("A627").Value = 10
$ws.Range("B627").Value = "Vega Modelo de Temuco"
$ws.Range("C627").Value = "La Araucanía"
$ws.Range("D627").Value = 45180
$ws.Range("E627").Value = 9
$ws.Range("F627").Value = 100112008
$ws.Range("G627").Value = "Coliflor"
$ws.Range("H627").Value = "Sin especificar"
$ws.Range("I627").Value = "Primera"
$ws.Range("J627").Value = 900
$ws.Range("K627").Value = 1500
$ws.Range("L627").Value = 1500
$ws.Range("M627").Value = 1500
$ws.Range("N627").Value = "`$/unidad"
$ws.Range("O627").Value = "Región Metropolitana"
$ws.Range("P627").Value = 1500
$ws.Range("Q627").Value = 1
$ws.Range("R627").Value = "Hortaliza"
